# Apply updated Betfair odds values to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "P2"  = 2.14
    "Q2"  = 1.82
    "R2"  = 1.45
    "S2"  = 3.05
    "T2"  = 1.69
    "U2"  = 2.34
    "N3"  = 1.99
    "Y4"  = 46
    "Z4"  = 80
    "G5"  = 2.46
    "J5"  = 3.1
    "W5"  = 1.68
    "Q6"  = 2.2
    "W6"  = 2.04
    "H7"  = 4.7
    "M7"  = 1.07
    "N7"  = 3.65
    "R7"  = 1.35
    "S7"  = 3.25
    "T7"  = 1.83
    "U7"  = 1.96
    "Z7"  = 980
    "AC7" = 1000
    "T8"  = 1.83
    "F9"  = 1.9
    "N9"  = 1.01
    "Q9"  = 1.38
    "V9"  = 1.22
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
